$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 / J1 ---------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold / bordered / centered)
# by copying the format from an existing header cell (H1) rather than
# re-building it from scratch, so the shared style id is reused.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data columns I2:J76 ------------------------------------------------
$iVals = @(9,9,10,8,9,9,9,8,9,9,8,9,9,9,9,8,9,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,8,9,9,9,10,9,9,10,8,9,9,9,9,9,9,9,9,9,9,8,9,9,9,9,9,9,8,9,8,9,9,9,10,9,9,9,9,9,9,9,9,5,3)
$jVals = @(9,9,10,9,9,9,10,9,9,9,8,9,9,9,9,9,10,9,9,9,9,9,9,9,9,10,9,9,9,9,9,9,9,9,9,9,10,9,9,10,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,9,8,9,9,10,11,9,9,9,9,9,9,9,9,5,3)

for ($r = 2; $r -le 76; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
